$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 69
$ws.Cells.Item(69, 2).Value = 3201281
$ws.Cells.Item(69, 6).Value = 'Coquimbo Unido'
$ws.Cells.Item(69, 7).Value = 'Everton de Vina'
$ws.Cells.Item(69, 8).Value = 3
$ws.Cells.Item(69, 9).Value = 0
$ws.Cells.Item(69, 10).Value = 'H'
$ws.Cells.Item(69, 11).Value = 3
$ws.Cells.Item(69, 12).Value = 3.2
$ws.Cells.Item(69, 13).Value = 2.375
$ws.Cells.Item(69, 14).Value = 3.5
$ws.Cells.Item(69, 15).Value = 3.2
$ws.Cells.Item(69, 16).Value = 2.15
$ws.Cells.Item(69, 17).Value = 0.25
$ws.Cells.Item(69, 18).Value = 1.925
$ws.Cells.Item(69, 19).Value = 1.875
$ws.Cells.Item(69, 20).Value = 2.25
$ws.Cells.Item(69, 21).Value = 1.8
$ws.Cells.Item(69, 22).Value = 2
$ws.Cells.Item(69, 23).Value = 2.5
$ws.Cells.Item(69, 24).Value = -1
$ws.Cells.Item(69, 25).Value = -1
$ws.Cells.Item(69, 26).Value = 0.925
$ws.Cells.Item(69, 27).Value = -1
$ws.Cells.Item(69, 28).Value = 0.8
$ws.Cells.Item(69, 29).Value = -1

# Row 70
$ws.Cells.Item(70, 2).Value = 3201280
$ws.Cells.Item(70, 6).Value = 'Santiago Wanderers'
$ws.Cells.Item(70, 7).Value = 'Universidad de Concepcion'
$ws.Cells.Item(70, 8).Value = 1
$ws.Cells.Item(70, 9).Value = 1
$ws.Cells.Item(70, 10).Value = 'D'
$ws.Cells.Item(70, 11).Value = 1.833
$ws.Cells.Item(70, 12).Value = 3.75
$ws.Cells.Item(70, 13).Value = 4
$ws.Cells.Item(70, 14).Value = 2
$ws.Cells.Item(70, 15).Value = 3.5
$ws.Cells.Item(70, 16).Value = 3.5
$ws.Cells.Item(70, 17).Value = -0.25
$ws.Cells.Item(70, 18).Value = 1.85
$ws.Cells.Item(70, 19).Value = 1.95
$ws.Cells.Item(70, 20).Value = 2.5
$ws.Cells.Item(70, 21).Value = 1.825
$ws.Cells.Item(70, 22).Value = 1.975
$ws.Cells.Item(70, 23).Value = -1
$ws.Cells.Item(70, 24).Value = 2.5
$ws.Cells.Item(70, 25).Value = -1
$ws.Cells.Item(70, 26).Value = -0.5
$ws.Cells.Item(70, 27).Value = 0.475
$ws.Cells.Item(70, 28).Value = -1
$ws.Cells.Item(70, 29).Value = 0.9750000000000001

# Row 71
$ws.Cells.Item(71, 2).Value = 3200931
$ws.Cells.Item(71, 6).Value = 'Colo Colo'
$ws.Cells.Item(71, 7).Value = 'Cobresal'
$ws.Cells.Item(71, 8).Value = 0
$ws.Cells.Item(71, 9).Value = 0
$ws.Cells.Item(71, 10).Value = 'D'
$ws.Cells.Item(71, 11).Value = 2.2
$ws.Cells.Item(71, 12).Value = 3.3
$ws.Cells.Item(71, 13).Value = 3.3
$ws.Cells.Item(71, 14).Value = 2
$ws.Cells.Item(71, 15).Value = 3.4
$ws.Cells.Item(71, 16).Value = 3.6
$ws.Cells.Item(71, 17).Value = -0.25
$ws.Cells.Item(71, 18).Value = 1.8
$ws.Cells.Item(71, 19).Value = 2
$ws.Cells.Item(71, 20).Value = 2.25
$ws.Cells.Item(71, 21).Value = 1.85
$ws.Cells.Item(71, 22).Value = 1.95
$ws.Cells.Item(71, 23).Value = -1
$ws.Cells.Item(71, 24).Value = 2.4
$ws.Cells.Item(71, 25).Value = -1
$ws.Cells.Item(71, 26).Value = -0.5
$ws.Cells.Item(71, 27).Value = 0.5
$ws.Cells.Item(71, 28).Value = -1
$ws.Cells.Item(71, 29).Value = 0.95

# Row 72
$ws.Cells.Item(72, 2).Value = 3201269
$ws.Cells.Item(72, 6).Value = 'Audax Italiano'
$ws.Cells.Item(72, 7).Value = 'Deportes Iquique'
$ws.Cells.Item(72, 8).Value = 1
$ws.Cells.Item(72, 9).Value = 1
$ws.Cells.Item(72, 10).Value = 'D'
$ws.Cells.Item(72, 11).Value = 2
$ws.Cells.Item(72, 12).Value = 3.5
$ws.Cells.Item(72, 13).Value = 3.6
$ws.Cells.Item(72, 14).Value = 2.3
$ws.Cells.Item(72, 15).Value = 3.5
$ws.Cells.Item(72, 16).Value = 2.875
$ws.Cells.Item(72, 17).Value = 0
$ws.Cells.Item(72, 18).Value = 1.75
$ws.Cells.Item(72, 19).Value = 2.05
$ws.Cells.Item(72, 20).Value = 2.5
$ws.Cells.Item(72, 21).Value = 1.825
$ws.Cells.Item(72, 22).Value = 1.975
$ws.Cells.Item(72, 23).Value = -1
$ws.Cells.Item(72, 24).Value = 2.5
$ws.Cells.Item(72, 25).Value = -1
$ws.Cells.Item(72, 26).Value = 0
$ws.Cells.Item(72, 27).Value = -0
$ws.Cells.Item(72, 28).Value = -1
$ws.Cells.Item(72, 29).Value = 0.9750000000000001

# Row 74
$ws.Cells.Item(74, 2).Value = 3234177
$ws.Cells.Item(74, 6).Value = 'Union La Calera'
$ws.Cells.Item(74, 7).Value = 'Curico Unido'
$ws.Cells.Item(74, 8).Value = 2
$ws.Cells.Item(74, 9).Value = 3
$ws.Cells.Item(74, 10).Value = 'A'
$ws.Cells.Item(74, 11).Value = 1.85
$ws.Cells.Item(74, 12).Value = 3.6
$ws.Cells.Item(74, 13).Value = 4
$ws.Cells.Item(74, 14).Value = 1.95
$ws.Cells.Item(74, 15).Value = 3.5
$ws.Cells.Item(74, 16).Value = 4
$ws.Cells.Item(74, 17).Value = -0.5
$ws.Cells.Item(74, 18).Value = 1.925
$ws.Cells.Item(74, 19).Value = 1.875
$ws.Cells.Item(74, 20).Value = 3
$ws.Cells.Item(74, 21).Value = 2
$ws.Cells.Item(74, 22).Value = 1.8
$ws.Cells.Item(74, 23).Value = -1
$ws.Cells.Item(74, 24).Value = -1
$ws.Cells.Item(74, 25).Value = 3
$ws.Cells.Item(74, 26).Value = -1
$ws.Cells.Item(74, 27).Value = 0.875
$ws.Cells.Item(74, 28).Value = 1
$ws.Cells.Item(74, 29).Value = -1

# Row 75
$ws.Cells.Item(75, 2).Value = 3234974
$ws.Cells.Item(75, 6).Value = 'Universidad de Chile'
$ws.Cells.Item(75, 7).Value = 'CD Antofagasta'
$ws.Cells.Item(75, 8).Value = 3
$ws.Cells.Item(75, 9).Value = 1
$ws.Cells.Item(75, 10).Value = 'H'
$ws.Cells.Item(75, 11).Value = 1.75
$ws.Cells.Item(75, 12).Value = 3.4
$ws.Cells.Item(75, 13).Value = 5
$ws.Cells.Item(75, 14).Value = 1.6
$ws.Cells.Item(75, 15).Value = 3.6
$ws.Cells.Item(75, 16).Value = 6
$ws.Cells.Item(75, 17).Value = -0.75
$ws.Cells.Item(75, 18).Value = 1.8
$ws.Cells.Item(75, 19).Value = 2
$ws.Cells.Item(75, 20).Value = 2.75
$ws.Cells.Item(75, 21).Value = 1.8
$ws.Cells.Item(75, 22).Value = 2
$ws.Cells.Item(75, 23).Value = 0.6000000000000001
$ws.Cells.Item(75, 24).Value = -1
$ws.Cells.Item(75, 25).Value = -1
$ws.Cells.Item(75, 26).Value = 0.8
$ws.Cells.Item(75, 27).Value = -1
$ws.Cells.Item(75, 28).Value = 0.8
$ws.Cells.Item(75, 29).Value = -1

# Row 346
$ws.Cells.Item(346, 2).Value = 3898955
$ws.Cells.Item(346, 6).Value = 'Universidad Catolica'
$ws.Cells.Item(346, 7).Value = 'Huachipato'
$ws.Cells.Item(346, 8).Value = 2
$ws.Cells.Item(346, 9).Value = 0
$ws.Cells.Item(346, 10).Value = 'H'
$ws.Cells.Item(346, 11).Value = 1.5
$ws.Cells.Item(346, 12).Value = 4
$ws.Cells.Item(346, 13).Value = 5.5
$ws.Cells.Item(346, 14).Value = 1.65
$ws.Cells.Item(346, 15).Value = 3.6
$ws.Cells.Item(346, 16).Value = 4.5
$ws.Cells.Item(346, 17).Value = -0.75
$ws.Cells.Item(346, 18).Value = 1.925
$ws.Cells.Item(346, 19).Value = 1.875
$ws.Cells.Item(346, 20).Value = 2.75
$ws.Cells.Item(346, 21).Value = 2
$ws.Cells.Item(346, 22).Value = 1.8
$ws.Cells.Item(346, 23).Value = 0.6499999999999999
$ws.Cells.Item(346, 24).Value = -1
$ws.Cells.Item(346, 25).Value = -1
$ws.Cells.Item(346, 26).Value = 0.925
$ws.Cells.Item(346, 27).Value = -1
$ws.Cells.Item(346, 28).Value = -1
$ws.Cells.Item(346, 29).Value = 0.8

# Row 347
$ws.Cells.Item(347, 2).Value = 4365259
$ws.Cells.Item(347, 6).Value = 'Curico Unido'
$ws.Cells.Item(347, 7).Value = 'Palestino'
$ws.Cells.Item(347, 8).Value = 1
$ws.Cells.Item(347, 9).Value = 1
$ws.Cells.Item(347, 10).Value = 'D'
$ws.Cells.Item(347, 11).Value = 2.45
$ws.Cells.Item(347, 12).Value = 3.25
$ws.Cells.Item(347, 13).Value = 2.625
$ws.Cells.Item(347, 14).Value = 2.25
$ws.Cells.Item(347, 15).Value = 3.25
$ws.Cells.Item(347, 16).Value = 2.9
$ws.Cells.Item(347, 17).Value = -0.25
$ws.Cells.Item(347, 18).Value = 1.95
$ws.Cells.Item(347, 19).Value = 1.85
$ws.Cells.Item(347, 20).Value = 2.5
$ws.Cells.Item(347, 21).Value = 1.925
$ws.Cells.Item(347, 22).Value = 1.875
$ws.Cells.Item(347, 23).Value = -1
$ws.Cells.Item(347, 24).Value = 2.25
$ws.Cells.Item(347, 25).Value = -1
$ws.Cells.Item(347, 26).Value = -0.5
$ws.Cells.Item(347, 27).Value = 0.425
$ws.Cells.Item(347, 28).Value = -1
$ws.Cells.Item(347, 29).Value = 0.875

# Row 348
$ws.Cells.Item(348, 2).Value = 3899776
$ws.Cells.Item(348, 6).Value = 'La Serena'
$ws.Cells.Item(348, 7).Value = 'Santiago Wanderers'
$ws.Cells.Item(348, 8).Value = 0
$ws.Cells.Item(348, 9).Value = 0
$ws.Cells.Item(348, 10).Value = 'D'
$ws.Cells.Item(348, 11).Value = 1.666
$ws.Cells.Item(348, 12).Value = 3.6
$ws.Cells.Item(348, 13).Value = 4.5
$ws.Cells.Item(348, 14).Value = 1.333
$ws.Cells.Item(348, 15).Value = 4.75
$ws.Cells.Item(348, 16).Value = 9.5
$ws.Cells.Item(348, 17).Value = -1.25
$ws.Cells.Item(348, 18).Value = 1.8
$ws.Cells.Item(348, 19).Value = 2
$ws.Cells.Item(348, 20).Value = 2.5
$ws.Cells.Item(348, 21).Value = 1.85
$ws.Cells.Item(348, 22).Value = 1.95
$ws.Cells.Item(348, 23).Value = -1
$ws.Cells.Item(348, 24).Value = 3.75
$ws.Cells.Item(348, 25).Value = -1
$ws.Cells.Item(348, 26).Value = -1
$ws.Cells.Item(348, 27).Value = 1
$ws.Cells.Item(348, 28).Value = -1
$ws.Cells.Item(348, 29).Value = 0.95

# Row 349
$ws.Cells.Item(349, 2).Value = 3898957
$ws.Cells.Item(349, 6).Value = 'Everton de Vina'
$ws.Cells.Item(349, 7).Value = 'Universidad Catolica'
$ws.Cells.Item(349, 8).Value = 0
$ws.Cells.Item(349, 9).Value = 3
$ws.Cells.Item(349, 10).Value = 'A'
$ws.Cells.Item(349, 11).Value = 3.25
$ws.Cells.Item(349, 12).Value = 3.25
$ws.Cells.Item(349, 13).Value = 2.05
$ws.Cells.Item(349, 14).Value = 3.5
$ws.Cells.Item(349, 15).Value = 3.4
$ws.Cells.Item(349, 16).Value = 2.1
$ws.Cells.Item(349, 17).Value = 0.25
$ws.Cells.Item(349, 18).Value = 2.05
$ws.Cells.Item(349, 19).Value = 1.8
$ws.Cells.Item(349, 20).Value = 2.25
$ws.Cells.Item(349, 21).Value = 1.825
$ws.Cells.Item(349, 22).Value = 2.025
$ws.Cells.Item(349, 23).Value = -1
$ws.Cells.Item(349, 24).Value = -1
$ws.Cells.Item(349, 25).Value = 1.1
$ws.Cells.Item(349, 26).Value = -1
$ws.Cells.Item(349, 27).Value = 0.8
$ws.Cells.Item(349, 28).Value = 0.825
$ws.Cells.Item(349, 29).Value = -1

# Row 350
$ws.Cells.Item(350, 2).Value = 3899774
$ws.Cells.Item(350, 6).Value = 'CD Antofagasta'
$ws.Cells.Item(350, 7).Value = 'Colo Colo'
$ws.Cells.Item(350, 8).Value = 1
$ws.Cells.Item(350, 9).Value = 0
$ws.Cells.Item(350, 10).Value = 'H'
$ws.Cells.Item(350, 11).Value = 4
$ws.Cells.Item(350, 12).Value = 3.4
$ws.Cells.Item(350, 13).Value = 1.8
$ws.Cells.Item(350, 14).Value = 9
$ws.Cells.Item(350, 15).Value = 4.75
$ws.Cells.Item(350, 16).Value = 1.363
$ws.Cells.Item(350, 17).Value = 1.5
$ws.Cells.Item(350, 18).Value = 1.775
$ws.Cells.Item(350, 19).Value = 2.025
$ws.Cells.Item(350, 20).Value = 2.5
$ws.Cells.Item(350, 21).Value = 1.8
$ws.Cells.Item(350, 22).Value = 2
$ws.Cells.Item(350, 23).Value = 8
$ws.Cells.Item(350, 24).Value = -1
$ws.Cells.Item(350, 25).Value = -1
$ws.Cells.Item(350, 26).Value = 0.7749999999999999
$ws.Cells.Item(350, 27).Value = -1
$ws.Cells.Item(350, 28).Value = -1
$ws.Cells.Item(350, 29).Value = 1

# Row 351
$ws.Cells.Item(351, 2).Value = 3899778
$ws.Cells.Item(351, 6).Value = 'Palestino'
$ws.Cells.Item(351, 7).Value = 'Nublense'
$ws.Cells.Item(351, 8).Value = 1
$ws.Cells.Item(351, 9).Value = 4
$ws.Cells.Item(351, 10).Value = 'A'
$ws.Cells.Item(351, 11).Value = 2.05
$ws.Cells.Item(351, 12).Value = 3.4
$ws.Cells.Item(351, 13).Value = 3.1
$ws.Cells.Item(351, 14).Value = 2
$ws.Cells.Item(351, 15).Value = 3.6
$ws.Cells.Item(351, 16).Value = 3.6
$ws.Cells.Item(351, 17).Value = -0.5
$ws.Cells.Item(351, 18).Value = 2
$ws.Cells.Item(351, 19).Value = 1.8
$ws.Cells.Item(351, 20).Value = 2.75
$ws.Cells.Item(351, 21).Value = 1.85
$ws.Cells.Item(351, 22).Value = 1.95
$ws.Cells.Item(351, 23).Value = -1
$ws.Cells.Item(351, 24).Value = -1
$ws.Cells.Item(351, 25).Value = 2.6
$ws.Cells.Item(351, 26).Value = -1
$ws.Cells.Item(351, 27).Value = 0.8
$ws.Cells.Item(351, 28).Value = 0.8500000000000001
$ws.Cells.Item(351, 29).Value = -1

# Row 587
$ws.Cells.Item(587, 2).Value = 4617491
$ws.Cells.Item(587, 6).Value = 'Universidad de Chile'
$ws.Cells.Item(587, 7).Value = 'Cobresal'
$ws.Cells.Item(587, 8).Value = 3
$ws.Cells.Item(587, 9).Value = 4
$ws.Cells.Item(587, 10).Value = 'A'
$ws.Cells.Item(587, 11).Value = 2.7
$ws.Cells.Item(587, 12).Value = 3.3
$ws.Cells.Item(587, 13).Value = 2.45
$ws.Cells.Item(587, 14).Value = 2.7
$ws.Cells.Item(587, 15).Value = 3.4
$ws.Cells.Item(587, 16).Value = 2.6
$ws.Cells.Item(587, 17).Value = 0
$ws.Cells.Item(587, 18).Value = 2
$ws.Cells.Item(587, 19).Value = 1.8
$ws.Cells.Item(587, 20).Value = 2.5
$ws.Cells.Item(587, 21).Value = 1.9
$ws.Cells.Item(587, 22).Value = 1.9
$ws.Cells.Item(587, 23).Value = -1
$ws.Cells.Item(587, 24).Value = -1
$ws.Cells.Item(587, 25).Value = 1.6
$ws.Cells.Item(587, 26).Value = -1
$ws.Cells.Item(587, 27).Value = 0.8
$ws.Cells.Item(587, 28).Value = 0.8999999999999999
$ws.Cells.Item(587, 29).Value = -1

# Row 588
$ws.Cells.Item(588, 2).Value = 4617747
$ws.Cells.Item(588, 6).Value = 'OHiggins'
$ws.Cells.Item(588, 7).Value = 'Everton de Vina'
$ws.Cells.Item(588, 8).Value = 2
$ws.Cells.Item(588, 9).Value = 0
$ws.Cells.Item(588, 10).Value = 'H'
$ws.Cells.Item(588, 11).Value = 2.4
$ws.Cells.Item(588, 12).Value = 3.1
$ws.Cells.Item(588, 13).Value = 2.9
$ws.Cells.Item(588, 14).Value = 2.15
$ws.Cells.Item(588, 15).Value = 3.25
$ws.Cells.Item(588, 16).Value = 3.5
$ws.Cells.Item(588, 17).Value = -0.25
$ws.Cells.Item(588, 18).Value = 1.825
$ws.Cells.Item(588, 19).Value = 1.975
$ws.Cells.Item(588, 20).Value = 2
$ws.Cells.Item(588, 21).Value = 1.775
$ws.Cells.Item(588, 22).Value = 2.025
$ws.Cells.Item(588, 23).Value = 1.15
$ws.Cells.Item(588, 24).Value = -1
$ws.Cells.Item(588, 25).Value = -1
$ws.Cells.Item(588, 26).Value = 0.825
$ws.Cells.Item(588, 27).Value = -1
$ws.Cells.Item(588, 28).Value = 0
$ws.Cells.Item(588, 29).Value = -0

# Row 589
$ws.Cells.Item(589, 2).Value = 4614405
$ws.Cells.Item(589, 6).Value = 'Palestino'
$ws.Cells.Item(589, 7).Value = 'Huachipato'
$ws.Cells.Item(589, 8).Value = 5
$ws.Cells.Item(589, 9).Value = 0
$ws.Cells.Item(589, 10).Value = 'H'
$ws.Cells.Item(589, 11).Value = 1.75
$ws.Cells.Item(589, 12).Value = 3.5
$ws.Cells.Item(589, 13).Value = 4.333
$ws.Cells.Item(589, 14).Value = 1.833
$ws.Cells.Item(589, 15).Value = 3.5
$ws.Cells.Item(589, 16).Value = 4.5
$ws.Cells.Item(589, 17).Value = -0.5
$ws.Cells.Item(589, 18).Value = 1.8
$ws.Cells.Item(589, 19).Value = 2
$ws.Cells.Item(589, 20).Value = 2.5
$ws.Cells.Item(589, 21).Value = 1.875
$ws.Cells.Item(589, 22).Value = 1.925
$ws.Cells.Item(589, 23).Value = 0.833
$ws.Cells.Item(589, 24).Value = -1
$ws.Cells.Item(589, 25).Value = -1
$ws.Cells.Item(589, 26).Value = 0.8
$ws.Cells.Item(589, 27).Value = -1
$ws.Cells.Item(589, 28).Value = 0.875
$ws.Cells.Item(589, 29).Value = -1

# Row 591
$ws.Cells.Item(591, 2).Value = 4614403
$ws.Cells.Item(591, 6).Value = 'CD Antofagasta'
$ws.Cells.Item(591, 7).Value = 'Universidad Catolica'
$ws.Cells.Item(591, 8).Value = 0
$ws.Cells.Item(591, 9).Value = 2
$ws.Cells.Item(591, 10).Value = 'A'
$ws.Cells.Item(591, 11).Value = 3
$ws.Cells.Item(591, 12).Value = 3.4
$ws.Cells.Item(591, 13).Value = 2.2
$ws.Cells.Item(591, 14).Value = 3.4
$ws.Cells.Item(591, 15).Value = 3.5
$ws.Cells.Item(591, 16).Value = 2.15
$ws.Cells.Item(591, 17).Value = 0.25
$ws.Cells.Item(591, 18).Value = 1.975
$ws.Cells.Item(591, 19).Value = 1.825
$ws.Cells.Item(591, 20).Value = 2.25
$ws.Cells.Item(591, 21).Value = 1.775
$ws.Cells.Item(591, 22).Value = 2.025
$ws.Cells.Item(591, 23).Value = -1
$ws.Cells.Item(591, 24).Value = -1
$ws.Cells.Item(591, 25).Value = 1.15
$ws.Cells.Item(591, 26).Value = -1
$ws.Cells.Item(591, 27).Value = 0.825
$ws.Cells.Item(591, 28).Value = -0.5
$ws.Cells.Item(591, 29).Value = 0.5125

# Row 592
$ws.Cells.Item(592, 2).Value = 4617739
$ws.Cells.Item(592, 6).Value = 'Coquimbo Unido'
$ws.Cells.Item(592, 7).Value = 'Curico Unido'
$ws.Cells.Item(592, 8).Value = 1
$ws.Cells.Item(592, 9).Value = 1
$ws.Cells.Item(592, 10).Value = 'D'
$ws.Cells.Item(592, 11).Value = 2.3
$ws.Cells.Item(592, 12).Value = 3.2
$ws.Cells.Item(592, 13).Value = 3
$ws.Cells.Item(592, 14).Value = 2.1
$ws.Cells.Item(592, 15).Value = 3.5
$ws.Cells.Item(592, 16).Value = 3.4
$ws.Cells.Item(592, 17).Value = -0.25
$ws.Cells.Item(592, 18).Value = 1.8
$ws.Cells.Item(592, 19).Value = 2
$ws.Cells.Item(592, 20).Value = 2.25
$ws.Cells.Item(592, 21).Value = 1.8
$ws.Cells.Item(592, 22).Value = 2
$ws.Cells.Item(592, 23).Value = -1
$ws.Cells.Item(592, 24).Value = 2.5
$ws.Cells.Item(592, 25).Value = -1
$ws.Cells.Item(592, 26).Value = -0.5
$ws.Cells.Item(592, 27).Value = 0.5
$ws.Cells.Item(592, 28).Value = -0.5
$ws.Cells.Item(592, 29).Value = 0.5

# Row 593
$ws.Cells.Item(593, 2).Value = 4617490
$ws.Cells.Item(593, 6).Value = 'Audax Italiano'
$ws.Cells.Item(593, 7).Value = 'La Serena'
$ws.Cells.Item(593, 8).Value = 3
$ws.Cells.Item(593, 9).Value = 0
$ws.Cells.Item(593, 10).Value = 'H'
$ws.Cells.Item(593, 11).Value = 1.55
$ws.Cells.Item(593, 12).Value = 4
$ws.Cells.Item(593, 13).Value = 5.5
$ws.Cells.Item(593, 14).Value = 1.4
$ws.Cells.Item(593, 15).Value = 4.75
$ws.Cells.Item(593, 16).Value = 7.5
$ws.Cells.Item(593, 17).Value = -1.25
$ws.Cells.Item(593, 18).Value = 1.85
$ws.Cells.Item(593, 19).Value = 1.95
$ws.Cells.Item(593, 20).Value = 3
$ws.Cells.Item(593, 21).Value = 2
$ws.Cells.Item(593, 22).Value = 1.8
$ws.Cells.Item(593, 23).Value = 0.3999999999999999
$ws.Cells.Item(593, 24).Value = -1
$ws.Cells.Item(593, 25).Value = -1
$ws.Cells.Item(593, 26).Value = 0.8500000000000001
$ws.Cells.Item(593, 27).Value = -1
$ws.Cells.Item(593, 28).Value = 0
$ws.Cells.Item(593, 29).Value = -0

# Row 805
$ws.Cells.Item(805, 2).Value = 7323253
$ws.Cells.Item(805, 6).Value = 'Union Espanola'
$ws.Cells.Item(805, 7).Value = 'OHiggins'
$ws.Cells.Item(805, 8).Value = 3
$ws.Cells.Item(805, 9).Value = 3
$ws.Cells.Item(805, 10).Value = 'D'
$ws.Cells.Item(805, 11).Value = 2
$ws.Cells.Item(805, 12).Value = 3.4
$ws.Cells.Item(805, 13).Value = 3.5
$ws.Cells.Item(805, 14).Value = 2.1
$ws.Cells.Item(805, 15).Value = 3.5
$ws.Cells.Item(805, 16).Value = 3.75
$ws.Cells.Item(805, 17).Value = -0.5
$ws.Cells.Item(805, 18).Value = 2.025
$ws.Cells.Item(805, 19).Value = 1.775
$ws.Cells.Item(805, 20).Value = 2.5
$ws.Cells.Item(805, 21).Value = 1.95
$ws.Cells.Item(805, 22).Value = 1.85
$ws.Cells.Item(805, 23).Value = -1
$ws.Cells.Item(805, 24).Value = 2.5
$ws.Cells.Item(805, 25).Value = -1
$ws.Cells.Item(805, 26).Value = -1
$ws.Cells.Item(805, 27).Value = 0.7749999999999999
$ws.Cells.Item(805, 28).Value = 0.95
$ws.Cells.Item(805, 29).Value = -1

# Row 806
$ws.Cells.Item(806, 2).Value = 7323186
$ws.Cells.Item(806, 6).Value = 'Coquimbo Unido'
$ws.Cells.Item(806, 7).Value = 'Deportes Copiapo'
$ws.Cells.Item(806, 8).Value = 1
$ws.Cells.Item(806, 9).Value = 0
$ws.Cells.Item(806, 10).Value = 'H'
$ws.Cells.Item(806, 11).Value = 2
$ws.Cells.Item(806, 12).Value = 3.4
$ws.Cells.Item(806, 13).Value = 3.5
$ws.Cells.Item(806, 14).Value = 1.727
$ws.Cells.Item(806, 15).Value = 3.8
$ws.Cells.Item(806, 16).Value = 4.75
$ws.Cells.Item(806, 17).Value = -0.75
$ws.Cells.Item(806, 18).Value = 1.9
$ws.Cells.Item(806, 19).Value = 1.9
$ws.Cells.Item(806, 20).Value = 2.75
$ws.Cells.Item(806, 21).Value = 1.85
$ws.Cells.Item(806, 22).Value = 1.95
$ws.Cells.Item(806, 23).Value = 0.7270000000000001
$ws.Cells.Item(806, 24).Value = -1
$ws.Cells.Item(806, 25).Value = -1
$ws.Cells.Item(806, 26).Value = 0.45
$ws.Cells.Item(806, 27).Value = -0.5
$ws.Cells.Item(806, 28).Value = -1
$ws.Cells.Item(806, 29).Value = 0.95

# Row 816
$ws.Cells.Item(816, 2).Value = 6078266
$ws.Cells.Item(816, 6).Value = 'Palestino'
$ws.Cells.Item(816, 7).Value = 'Curico Unido'
$ws.Cells.Item(816, 8).Value = 4
$ws.Cells.Item(816, 9).Value = 0
$ws.Cells.Item(816, 10).Value = 'H'
$ws.Cells.Item(816, 11).Value = 1.533
$ws.Cells.Item(816, 12).Value = 4
$ws.Cells.Item(816, 13).Value = 6
$ws.Cells.Item(816, 14).Value = 1.363
$ws.Cells.Item(816, 15).Value = 4.75
$ws.Cells.Item(816, 16).Value = 7.5
$ws.Cells.Item(816, 17).Value = -1.5
$ws.Cells.Item(816, 18).Value = 2.025
$ws.Cells.Item(816, 19).Value = 1.825
$ws.Cells.Item(816, 20).Value = 3
$ws.Cells.Item(816, 21).Value = 1.9
$ws.Cells.Item(816, 22).Value = 1.95
$ws.Cells.Item(816, 23).Value = 0.363
$ws.Cells.Item(816, 24).Value = -1
$ws.Cells.Item(816, 25).Value = -1
$ws.Cells.Item(816, 26).Value = 1.025
$ws.Cells.Item(816, 27).Value = -1
$ws.Cells.Item(816, 28).Value = 0.8999999999999999
$ws.Cells.Item(816, 29).Value = -1

# Row 817
$ws.Cells.Item(817, 2).Value = 6077498
$ws.Cells.Item(817, 6).Value = 'Universidad Catolica'
$ws.Cells.Item(817, 7).Value = 'Deportes Copiapo'
$ws.Cells.Item(817, 8).Value = 2
$ws.Cells.Item(817, 9).Value = 2
$ws.Cells.Item(817, 10).Value = 'D'
$ws.Cells.Item(817, 11).Value = 1.65
$ws.Cells.Item(817, 12).Value = 3.8
$ws.Cells.Item(817, 13).Value = 5.25
$ws.Cells.Item(817, 14).Value = 1.909
$ws.Cells.Item(817, 15).Value = 3.6
$ws.Cells.Item(817, 16).Value = 4.2
$ws.Cells.Item(817, 17).Value = -0.5
$ws.Cells.Item(817, 18).Value = 1.85
$ws.Cells.Item(817, 19).Value = 2
$ws.Cells.Item(817, 20).Value = 2.75
$ws.Cells.Item(817, 21).Value = 2.025
$ws.Cells.Item(817, 22).Value = 1.825
$ws.Cells.Item(817, 23).Value = -1
$ws.Cells.Item(817, 24).Value = 2.6
$ws.Cells.Item(817, 25).Value = -1
$ws.Cells.Item(817, 26).Value = -1
$ws.Cells.Item(817, 27).Value = 1
$ws.Cells.Item(817, 28).Value = 1.025
$ws.Cells.Item(817, 29).Value = -1

# Row 818
$ws.Cells.Item(818, 2).Value = 6078265
$ws.Cells.Item(818, 6).Value = 'Audax Italiano'
$ws.Cells.Item(818, 7).Value = 'Magallanes'
$ws.Cells.Item(818, 8).Value = 0
$ws.Cells.Item(818, 9).Value = 2
$ws.Cells.Item(818, 10).Value = 'A'
$ws.Cells.Item(818, 11).Value = 1.666
$ws.Cells.Item(818, 12).Value = 3.75
$ws.Cells.Item(818, 13).Value = 5
$ws.Cells.Item(818, 14).Value = 2.25
$ws.Cells.Item(818, 15).Value = 3.3
$ws.Cells.Item(818, 16).Value = 3.3
$ws.Cells.Item(818, 17).Value = -0.25
$ws.Cells.Item(818, 18).Value = 1.95
$ws.Cells.Item(818, 19).Value = 1.85
$ws.Cells.Item(818, 20).Value = 2.5
$ws.Cells.Item(818, 21).Value = 1.8
$ws.Cells.Item(818, 22).Value = 2
$ws.Cells.Item(818, 23).Value = -1
$ws.Cells.Item(818, 24).Value = -1
$ws.Cells.Item(818, 25).Value = 2.3
$ws.Cells.Item(818, 26).Value = -1
$ws.Cells.Item(818, 27).Value = 0.8500000000000001
$ws.Cells.Item(818, 28).Value = -1
$ws.Cells.Item(818, 29).Value = 1

# Row 826
$ws.Cells.Item(826, 2).Value = 6078269
$ws.Cells.Item(826, 6).Value = 'Universidad de Chile'
$ws.Cells.Item(826, 7).Value = 'Nublense'
$ws.Cells.Item(826, 8).Value = 3
$ws.Cells.Item(826, 9).Value = 1
$ws.Cells.Item(826, 10).Value = 'H'
$ws.Cells.Item(826, 11).Value = 1.85
$ws.Cells.Item(826, 12).Value = 3.4
$ws.Cells.Item(826, 13).Value = 4.333
$ws.Cells.Item(826, 14).Value = 1.8
$ws.Cells.Item(826, 15).Value = 3.6
$ws.Cells.Item(826, 16).Value = 4.5
$ws.Cells.Item(826, 17).Value = -0.75
$ws.Cells.Item(826, 18).Value = 1.925
$ws.Cells.Item(826, 19).Value = 1.925
$ws.Cells.Item(826, 20).Value = 2.5
$ws.Cells.Item(826, 21).Value = 2.025
$ws.Cells.Item(826, 22).Value = 1.825
$ws.Cells.Item(826, 23).Value = 0.8
$ws.Cells.Item(826, 24).Value = -1
$ws.Cells.Item(826, 25).Value = -1
$ws.Cells.Item(826, 26).Value = 0.925
$ws.Cells.Item(826, 27).Value = -1
$ws.Cells.Item(826, 28).Value = 1.025
$ws.Cells.Item(826, 29).Value = -1

# Row 827
$ws.Cells.Item(827, 2).Value = 6078998
$ws.Cells.Item(827, 6).Value = 'Magallanes'
$ws.Cells.Item(827, 7).Value = 'Coquimbo Unido'
$ws.Cells.Item(827, 8).Value = 2
$ws.Cells.Item(827, 9).Value = 3
$ws.Cells.Item(827, 10).Value = 'A'
$ws.Cells.Item(827, 11).Value = 1.909
$ws.Cells.Item(827, 12).Value = 3.6
$ws.Cells.Item(827, 13).Value = 3.8
$ws.Cells.Item(827, 14).Value = 2.15
$ws.Cells.Item(827, 15).Value = 3.75
$ws.Cells.Item(827, 16).Value = 3.1
$ws.Cells.Item(827, 17).Value = -0.25
$ws.Cells.Item(827, 18).Value = 1.85
$ws.Cells.Item(827, 19).Value = 1.95
$ws.Cells.Item(827, 20).Value = 3
$ws.Cells.Item(827, 21).Value = 1.85
$ws.Cells.Item(827, 22).Value = 1.95
$ws.Cells.Item(827, 23).Value = -1
$ws.Cells.Item(827, 24).Value = -1
$ws.Cells.Item(827, 25).Value = 2.1
$ws.Cells.Item(827, 26).Value = -1
$ws.Cells.Item(827, 27).Value = 0.95
$ws.Cells.Item(827, 28).Value = 0.8500000000000001
$ws.Cells.Item(827, 29).Value = -1

# Row 828
$ws.Cells.Item(828, 2).Value = 6077499
$ws.Cells.Item(828, 6).Value = 'Deportes Copiapo'
$ws.Cells.Item(828, 7).Value = 'Everton de Vina'
$ws.Cells.Item(828, 8).Value = 2
$ws.Cells.Item(828, 9).Value = 0
$ws.Cells.Item(828, 10).Value = 'H'
$ws.Cells.Item(828, 11).Value = 2.1
$ws.Cells.Item(828, 12).Value = 3.4
$ws.Cells.Item(828, 13).Value = 3.4
$ws.Cells.Item(828, 14).Value = 2.9
$ws.Cells.Item(828, 15).Value = 3.4
$ws.Cells.Item(828, 16).Value = 2.4
$ws.Cells.Item(828, 17).Value = 0.25
$ws.Cells.Item(828, 18).Value = 1.775
$ws.Cells.Item(828, 19).Value = 2.1
$ws.Cells.Item(828, 20).Value = 2.75
$ws.Cells.Item(828, 21).Value = 1.85
$ws.Cells.Item(828, 22).Value = 2
$ws.Cells.Item(828, 23).Value = 1.9
$ws.Cells.Item(828, 24).Value = -1
$ws.Cells.Item(828, 25).Value = -1
$ws.Cells.Item(828, 26).Value = 0.7749999999999999
$ws.Cells.Item(828, 27).Value = -1
$ws.Cells.Item(828, 28).Value = -1
$ws.Cells.Item(828, 29).Value = 1

# Row 829
$ws.Cells.Item(829, 2).Value = 6078268
$ws.Cells.Item(829, 6).Value = 'OHiggins'
$ws.Cells.Item(829, 7).Value = 'Palestino'
$ws.Cells.Item(829, 8).Value = 0
$ws.Cells.Item(829, 9).Value = 1
$ws.Cells.Item(829, 10).Value = 'A'
$ws.Cells.Item(829, 11).Value = 3.1
$ws.Cells.Item(829, 12).Value = 3.3
$ws.Cells.Item(829, 13).Value = 2.3
$ws.Cells.Item(829, 14).Value = 2.9
$ws.Cells.Item(829, 15).Value = 3.4
$ws.Cells.Item(829, 16).Value = 2.375
$ws.Cells.Item(829, 17).Value = 0.25
$ws.Cells.Item(829, 18).Value = 1.8
$ws.Cells.Item(829, 19).Value = 2
$ws.Cells.Item(829, 20).Value = 2.75
$ws.Cells.Item(829, 21).Value = 2
$ws.Cells.Item(829, 22).Value = 1.8
$ws.Cells.Item(829, 23).Value = -1
$ws.Cells.Item(829, 24).Value = -1
$ws.Cells.Item(829, 25).Value = 1.375
$ws.Cells.Item(829, 26).Value = -1
$ws.Cells.Item(829, 27).Value = 1
$ws.Cells.Item(829, 28).Value = -1
$ws.Cells.Item(829, 29).Value = 0.8

# Row 830
$ws.Cells.Item(830, 2).Value = 6077768
$ws.Cells.Item(830, 6).Value = 'Union La Calera'
$ws.Cells.Item(830, 7).Value = 'Universidad Catolica'
$ws.Cells.Item(830, 8).Value = 0
$ws.Cells.Item(830, 9).Value = 3
$ws.Cells.Item(830, 10).Value = 'A'
$ws.Cells.Item(830, 11).Value = 2.05
$ws.Cells.Item(830, 12).Value = 3.5
$ws.Cells.Item(830, 13).Value = 3.4
$ws.Cells.Item(830, 14).Value = 2.05
$ws.Cells.Item(830, 15).Value = 3.6
$ws.Cells.Item(830, 16).Value = 3.4
$ws.Cells.Item(830, 17).Value = -0.25
$ws.Cells.Item(830, 18).Value = 1.8
$ws.Cells.Item(830, 19).Value = 2
$ws.Cells.Item(830, 20).Value = 2.75
$ws.Cells.Item(830, 21).Value = 1.975
$ws.Cells.Item(830, 22).Value = 1.825
$ws.Cells.Item(830, 23).Value = -1
$ws.Cells.Item(830, 24).Value = -1
$ws.Cells.Item(830, 25).Value = 2.4
$ws.Cells.Item(830, 26).Value = -1
$ws.Cells.Item(830, 27).Value = 1
$ws.Cells.Item(830, 28).Value = 0.4875
$ws.Cells.Item(830, 29).Value = -0.5
